$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; this shifts the existing
# B:F (PercActivations...) columns to C:G, carrying their styles
# and values along with them.
$ws.Columns("B:B").Insert()

# New column header, styled like the other header cells
# (bold, thin border all around, centered / top-aligned).
$hdr = $ws.Range("B1")
$hdr.Value = "segments"
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# Segment names, in row order (was column A text, now column B,
# alongside the new 0-based numeric index in column A).
$segments = @(
  "background",
  "back_bumper",
  "back_glass",
  "back_left_door",
  "back_left_light",
  "back_right_door",
  "back_right_light",
  "front_bumper",
  "front_glass",
  "front_left_door",
  "front_left_light",
  "front_right_door",
  "front_right_light",
  "hood",
  "left_mirror",
  "right_mirror",
  "tailgate",
  "trunk",
  "wheel"
)

for ($i = 0; $i -lt $segments.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $i
  $cell = $ws.Cells.Item($row, 2)
  $cell.Value = $segments[$i]
  # Plain formatting (matches the rest of the untouched data cells).
  $cell.Style = "Normal"
}
